# Adding new iAuthor testcases changes
#
# Replaces the small 5-row sample data block (rows 2-6) with the full
# 14-row iAuthor testcase export (rows 2-15), covering VSAQ / ISAWE /
# Type X / Type B / SAQ / SJT question types.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @(1,  1, "MCQ",    "Nill", "Nill", "Nill", "C",    8.5,    -0.06, 27.8,   -0.03, 25.1,  0.21,  21.3,  -0.07, 17.3,  -0.09, "Nill", "Nill", 0.21, 2.25, "Low")
  3  = @(2,  2, "VSAQ",   "Nill", "Nill", "Nill", "Nill", "Nill", "Nill","Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill", 55.0,  "Nill", 0.45, -0.9, "Low")
  4  = @(3,  3, "ISAWE",  1,      1,      "Nill", "Nill", "Nill", "Nill","Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill", 35.67, "Nill", 0.32, -0.9, "Low")
  5  = @(4,  3, "ISAWE",  1,      2,      "Nill", "Nill", "Nill", "Nill","Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill", 65.32, "Nill", 0.21, -0.9, "Low")
  6  = @(5,  3, "ISAWE",  2,      1,      "Nill", "Nill", "Nill", "Nill","Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill", 54.12, "Nill", 0.37, -0.9, "High")
  7  = @(6,  3, "ISAWE",  2,      2,      "Nill", "Nill", "Nill", "Nill","Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill", 38.05, "Nill", 0.54, -0.9, "High")
  8  = @(7,  4, "Type X", "Nill", "Nill", "A",    "F",    59.9,   -0.25, "Nill","Nill","Nill","Nill","Nill","Nill","Nill","Nill", "Nill","Nill", 0.26, 1.54, "Medium")
  9  = @(8,  4, "Type X", "Nill", "Nill", "B",    "T",    "Nill", "Nill",12.6,  -0.28, "Nill","Nill","Nill","Nill","Nill","Nill", "Nill","Nill", 0.28, -0.9, "Medium")
  10 = @(9,  4, "Type X", "Nill", "Nill", "C",    "T",    "Nill", "Nill","Nill","Nill",39.8,  0.26,  "Nill","Nill","Nill","Nill", "Nill","Nill", 0.28, -0.9, "Low")
  11 = @(10, 4, "Type X", "Nill", "Nill", "D",    "T",    "Nill", "Nill","Nill","Nill","Nill","Nill",12.6,  -0.28, "Nill","Nill", "Nill","Nill", 0.28, -0.9, "High")
  12 = @(11, 4, "Type X", "Nill", "Nill", "E",    "T",    "Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill",87.4,  0.28,  "Nill","Nill", 0.28, -0.9, "Low")
  13 = @(12, 5, "Type B", "Nill", "Nill", "Nill", "C",    8.5,    -0.06, 27.8,  -0.03, 25.1,  0.21,  21.3,  -0.07, 17.3,  -0.09, "Nill", "Nill", 0.21, 2.25, "Low")
  14 = @(13, 6, "SAQ",    "Nill", "Nill", "Nill", "Nill", "Nill", "Nill","Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill", 55.0,  "Nill", 0.45, -0.9, "Low")
  15 = @(14, 7, "SJT",    "Nill", "Nill", "Nill", "Nill", "Nill", "Nill","Nill", "Nill","Nill","Nill","Nill","Nill","Nill","Nill", 55.0,  "Nill", 0.45, -0.9, "Low")
}

foreach ($rowNum in 2..15) {
  $values = $data[$rowNum]
  for ($col = 1; $col -le $values.Length; $col++) {
    $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
  }
}
